# Update translation file: reorder/update a handful of rows in the
# en/kh/status translation sheet, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two obsolete "ACORN data is not of the right format..." rows
#    from their current position (row 11 twice, since after the first
#    delete the second target row becomes row 11).
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(11).Delete()

# 2) Insert a new row for "Contains names of organisms before and after
#    mapping." right before "Couldn't connect to server..." (now row 38).
$ws.Rows.Item(38).Insert()
$ws.Range("A38").Value = "Contains names of organisms before and after mapping."
$ws.Range("B38").Value = "TBT"
$ws.Range("C38").Value = "new"

# 3) Insert a new row for "Download Lab Log (.xlsx)" right after
#    "Download Enrolment Log (.xlsx)" (now row 50).
$ws.Rows.Item(50).Insert()
$ws.Range("A50").Value = "Download Lab Log (.xlsx)"
$ws.Range("B50").Value = "TBT"
$ws.Range("C50").Value = "new"

# 4) Replace row 70's content in place (old key retired, new key added).
$ws.Range("A70").Value = "HAI point prevalence by "
$ws.Range("B70").Value = "TBT"
$ws.Range("C70").Value = "new"

# 5) Insert a new row for "Remove 'Not Cultured' specimens" right before
#    "Remove blood culture contaminants..." (now row 110).
$ws.Rows.Item(110).Insert()
$ws.Range("A110").Value = "Remove 'Not Cultured' specimens"
$ws.Range("B110").Value = "TBT"
$ws.Range("C110").Value = "new"

# 6) Remove the obsolete "Select lab data format:" row (now row 125).
$ws.Rows.Item(125).Delete()

# 7) Append the 4 retired rows at the bottom (182-185) marked "deleted",
#    keeping their original English/Khmer text.
$ws.Range("A182").Value = "ACORN data is not of the right format. Only data generated with v2.1 (or later versions) is compatible."
$ws.Range("B182").Value = "TBT"
$ws.Range("C182").Value = "deleted"

$ws.Range("A183").Value = "ACORN data is not of the right format. Only data generated with v2.1 is compatible."
$ws.Range("B183").Value = "ទិន្នន័យ ACORN មិនមែនជាទម្រង់ត្រឹមត្រូវទេ។ មានតែទិន្នន័យដែលបង្កើតឡើងដោយជំនាន់ V២.១ ប៉ុណ្ណោះដែលត្រូវគ្នា។"
$ws.Range("C183").Value = "deleted"

$ws.Range("A184").Value = "HAI point prevalence by type of ward"
$ws.Range("B184").Value = "HAI point prevalance តាមប្រភេទអគារ"
$ws.Range("C184").Value = "deleted"

$ws.Range("A185").Value = "Select lab data format:"
$ws.Range("B185").Value = "ជ្រើសរើសទម្រង់នៃទិន្នន័យមន្ទីរពិសោធន៍៖"
$ws.Range("C185").Value = "deleted"
